$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 23:52"

# Update Cataluña row (row 5) figures
$ws.Range("B5").Value = 29647
$ws.Range("C5").Value = 12250
$ws.Range("D5").Value = 14356
$ws.Range("E5").Value = 3041
